$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 387, shifting existing rows 387:415 down to 388:416
$ws.Rows(387).Insert()

# Populate the newly inserted row 387 with the new data point
$ws.Range("A387").Value = 4
$ws.Range("B387").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C387").Value = "Los Lagos"
$ws.Range("D387").Value = 44783
$ws.Range("E387").Value = 10
$ws.Range("F387").Value = 100114013
$ws.Range("G387").Value = "Zanahoria"
$ws.Range("H387").Value = "Sin especificar"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 70
$ws.Range("K387").Value = 10000
$ws.Range("L387").Value = 10000
$ws.Range("M387").Value = 10000
$ws.Range("N387").Value = "`$/saco 20 kilos"
$ws.Range("O387").Value = "Provincia de Llanquihue"
$ws.Range("P387").Value = 500
$ws.Range("Q387").Value = 20
$ws.Range("R387").Value = "Hortaliza"
